# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 114
    4  = 1621
    8  = 11512
    11 = 450
    12 = 362
    14 = 802
    16 = 13084
    21 = 226
    24 = 119
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
